$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.736.26"
$ws.Cells.Item(2, 5).Value = "  -3.09%  "

$ws.Cells.Item(3, 4).Value = "2.609.76"
$ws.Cells.Item(3, 5).Value = "  -2.13%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).Value = "'574.28"
$ws.Cells.Item(5, 5).Value = "  -3.82%  "

$ws.Cells.Item(6, 4).Value = "'155.75"
$ws.Cells.Item(6, 5).Value = "  -0.95%  "

$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  +0.00%  "

$ws.Cells.Item(8, 4).Value = "'0.620"
$ws.Cells.Item(8, 5).Value = "  -6.13%  "

$ws.Cells.Item(9, 5).Value = "  -5.95%  "

$ws.Cells.Item(10, 5).Value = "  -0.14%  "

$ws.Cells.Item(11, 4).Value = "'0.381"
$ws.Cells.Item(11, 5).Value = "  -4.70%  "

$ws.Cells.Item(12, 5).Value = "  -0.31%  "

$ws.Cells.Item(13, 4).Value = "'28.09"
$ws.Cells.Item(13, 5).Value = "  -2.42%  "

$ws.Cells.Item(14, 4).Value = "3.077.62"
$ws.Cells.Item(14, 5).Value = "  -2.01%  "

$ws.Cells.Item(15, 4).Value = "'0.0000181"
$ws.Cells.Item(15, 5).Value = "  -6.92%  "

$ws.Cells.Item(16, 4).Value = "63.519.14"
$ws.Cells.Item(16, 5).Value = "  -3.20%  "

$ws.Cells.Item(17, 4).Value = "2.606.20"
$ws.Cells.Item(17, 5).Value = "  +0.20%  "

$ws.Cells.Item(18, 4).Value = "'12.03"
$ws.Cells.Item(18, 5).Value = "  -4.26%  "

$ws.Cells.Item(19, 4).Value = "'7.62"
$ws.Cells.Item(19, 5).Value = "  +2.64%  "

$ws.Cells.Item(20, 4).Value = "'4.54"
$ws.Cells.Item(20, 5).Value = "  -5.02%  "

$ws.Cells.Item(21, 4).Value = "'342.36"
$ws.Cells.Item(21, 5).Value = "  -2.10%  "

$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 5).Value = "  -0.10%  "

$ws.Cells.Item(23, 4).Value = "'67.12"
$ws.Cells.Item(23, 5).Value = "  -3.63%  "

$ws.Cells.Item(25, 4).Value = "'0.0000108"
$ws.Cells.Item(25, 5).Value = "  -3.19%  "

$ws.Cells.Item(26, 4).Value = "'588.74"
$ws.Cells.Item(26, 5).Value = "  +4.40%  "

$ws.Cells.Item(27, 4).Value = "'9.16"
$ws.Cells.Item(27, 5).Value = "  -3.89%  "

$ws.Cells.Item(28, 5).Value = "  -2.53%  "

$ws.Cells.Item(29, 5).Value = "  +0.16%  "

$ws.Cells.Item(30, 5).Value = "  -1.58%  "

$ws.Cells.Item(31, 4).Value = "'7.90"
$ws.Cells.Item(31, 5).Value = "  -1.71%  "

$ws.Cells.Item(32, 5).Value = "  -4.68%  "

$ws.Cells.Item(33, 4).Value = "'1.74"
$ws.Cells.Item(33, 5).Value = "  -3.57%  "

$ws.Cells.Item(34, 4).Value = "'6.52"
$ws.Cells.Item(34, 5).Value = "  -0.54%  "

$ws.Cells.Item(35, 5).Value = "  -2.27%  "

$ws.Cells.Item(36, 4).Value = "'0.406"
$ws.Cells.Item(36, 5).Value = "  -3.54%  "

$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 5).Value = "  -0.06%  "

$ws.Cells.Item(38, 4).Value = "'19.63"
$ws.Cells.Item(38, 5).Value = "  -4.13%  "

$ws.Cells.Item(39, 4).Value = "'154.05"
$ws.Cells.Item(39, 5).Value = "  -0.32%  "

$ws.Cells.Item(40, 5).Value = "  -3.79%  "

$ws.Cells.Item(41, 5).Value = "  -0.06%  "

$ws.Cells.Item(42, 4).Value = "'41.45"
$ws.Cells.Item(42, 5).Value = "  -2.97%  "

$ws.Cells.Item(43, 5).Value = "  +5.95%  "

$ws.Cells.Item(44, 4).Value = "'155.72"
$ws.Cells.Item(44, 5).Value = "  -3.07%  "

$ws.Cells.Item(45, 5).Value = "  -4.39%  "

$ws.Cells.Item(46, 4).Value = "'23.23"
$ws.Cells.Item(46, 5).Value = "  +2.49%  "

$ws.Cells.Item(47, 4).Value = "'0.0589"
$ws.Cells.Item(47, 5).Value = "  -2.46%  "

$ws.Cells.Item(48, 5).Value = "  -1.92%  "

$ws.Cells.Item(49, 4).Value = "'0.101"
$ws.Cells.Item(49, 5).Value = "  -1.71%  "

$ws.Cells.Item(50, 5).Value = "  -3.53%  "

$ws.Cells.Item(51, 4).Value = "'18.86"
$ws.Cells.Item(51, 5).Value = "  -4.66%  "
